# Updated cryptos list (price + volume refresh, and a rank swap for
# FTXToken / InjectiveProtocol) per the Fri Nov 24 20:44:00 UTC 2023 run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.859.87"
$ws.Range("E2").Value = "  +1.44%  "

$ws.Range("D3").Value = "2.086.75"
$ws.Range("E3").Value = "  +1.09%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.624"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.34"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.91%  "

$ws.Range("E9").Value = "  +1.72%  "

$ws.Range("E10").Value = "  +2.30%  "

$ws.Range("E11").Value = "  +2.92%  "

$ws.Range("D12").Value = "2.381.25"
$ws.Range("E12").Value = "  +0.57%  "

$ws.Range("E13").Value = "  -1.43%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.761"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.78%  "

$ws.Range("E16").Value = "  +2.04%  "

$ws.Range("D17").Value = "2.093.74"
$ws.Range("E17").Value = "  +1.42%  "

$ws.Range("D18").Value = "37.776.47"
$ws.Range("E18").Value = "  +1.36%  "

$ws.Range("E19").Value = "  -1.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.75"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").Value = "0.0₃0820"
$ws.Range("E21").Value = "  +1.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.88%  "

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("E24").Value = "  -1.66%  "

$ws.Range("E25").Value = "  -0.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.139"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.94%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.38%  "

$ws.Range("E29").Value = "  +0.67%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.20%  "

$ws.Range("E31").Value = "  +0.98%  "

$ws.Range("E32").Value = "  +3.97%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0626"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.88%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.10%  "

$ws.Range("E35").Value = "  +0.86%  "

$ws.Range("E36").Value = "  +3.63%  "

$ws.Range("E37").Value = "  +4.74%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0998"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.38%  "

$ws.Range("E41").Value = "  -0.87%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.09%  "

$ws.Range("E43").Value = "  +0.74%  "

$ws.Range("D44").Value = "1.451.11"
$ws.Range("E44").Value = "  -0.93%  "

$ws.Range("E45").Value = "  -0.22%  "

$ws.Range("E46").Value = "  +3.16%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.32%  "

$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.50%  "

$ws.Range("E50").Value = "  +1.51%  "

$ws.Range("D51").Value = "2.276.54"
$ws.Range("E51").Value = "  +0.85%  "
